$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.915.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.04%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.874.23"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.7412"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.15%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'242.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.61%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3143"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.53%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07167"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.96%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'24.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.73%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.08413"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.76%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.7526"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.48%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.400"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.13%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'1.877.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -9.26%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'92.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.86%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'29.917.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.61%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'6.100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.93%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -2.61%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'243.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.89%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.66%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.9990"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.16%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'2.125.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -8.25%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.68%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.9969"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.39%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -2.60%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'9.318"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.28%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'165.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.86%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -1.11%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -0.22%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.482"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.32%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'4.616"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.10%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.535"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.58%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'4.253"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.15%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.05324"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.64%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -0.73%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.7558"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.34%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.9972"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.77%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'2.697"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.26%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -1.09%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.749"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.36%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.4489"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.54%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'1.111.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.91%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'6.079"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.13%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'72.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.37%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.8581"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.34%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.07%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'103.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.12%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'7.665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.57%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'3.075"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.14%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.841"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.37%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'2.023.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -7.77%  "
$ws.Range("E51").Style = "Normal"
